# Add a new second slide ("How it works") to the deck, using the same
# "Title and Content" layout as the existing slide (slideLayout2.xml,
# the 2nd layout registered on the slide master).
$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(2, 2)

# Title placeholder
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "How it works"

# Body / content placeholder - four bullet paragraphs, 18pt like slide 1
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Stage metering " + [char]0x2192 + " TECHUP.AUDIT.WAREHOUSE_METERING_STG (task)`rRIGHT_SIZING_POLICY_DT: avg(credits_used) per warehouse/hour`rAPPLY_RIGHT_SIZING(): ALTER WAREHOUSE size + optional multi-cluster`rAPPLY_RIGHT_SIZING_TASK: runs hourly; logs actions to RIGHT_SIZING_LOG"
$body.Font.Size = 18
